# Apply the edits captured by the diff:
#  - update the USD Amount value in T2 (384005 -> 390200)
#  - move the active cell / selection to T3 (was R11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in cell T2
$ws.Range("T2").Value = 390200

# Update the current selection / active cell to T3
$ws.Range("T3").Select()
